# Adding maven surefire plugin
# -----------------------------------------------------------------------
# This script recreates, via Excel COM automation, the edit described by
# the target diff:
#   * Sheet1 header row becomes bold; the 3rd header cell's hyperlink is
#     removed but keeps its bold/underline/blue styling; many filler
#     header cells (D1:Y1) get the bold styling too (Sheet1 was widened).
#   * Row 2 values are swapped for "shivam"/"MyTest"/<MyTest url>.
#   * Row 2's extra columns D2/E2 are removed.
#   * A new Row 3 is added ("sachin"/"MyTest"/<MyTest url>), mirroring
#     row 2's hyperlink.
#   * Column C is widened.
#   * A new worksheet "Sheet2" is appended after Sheet1 with a small
#     "Test Case Repo" / "Test Case Repo URL" table.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Sheet1 - header row (row 1)
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "User Name"
$ws1.Range("B1").Value = "User Repo"
$ws1.Range("C1").Value = "User Repo URL"

# Headers become bold.
$ws1.Range("A1:C1").Font.Bold = $true

# C1 no longer is an actual hyperlink, but visually keeps the
# underline / link-blue color it already had.
$ws1.Range("C1").Hyperlinks.Delete()
$ws1.Range("C1").Font.Underline = $true
$ws1.Range("C1").Font.Color = 13391121

# D1/E1 lose their old values & hyperlinks (sheet now "ends" after C
# except for a band of blank, bold, themed filler cells D1:Y1).
$ws1.Range("D1:E1").Hyperlinks.Delete()
$ws1.Range("D1:Y1").ClearContents()
$ws1.Range("D1:Y1").Font.Bold = $true
$ws1.Range("D1:Y1").Font.Name = "Arial"

# ---------------------------------------------------------------------
# Sheet1 - data row 2
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "shivam"
$ws1.Range("B2").Value = "MyTest"
$ws1.Range("C2").Value = "https://github.com/shivamgupta2607/MyTest"

$ws1.Range("A2").Font.Name = "Arial"
$ws1.Range("A2:B2").Font.Bold = $false

$ws1.Range("C2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("C2"), "https://github.com/shivamgupta2607/MyTest") | Out-Null
$ws1.Range("C2").Font.Underline = $true
$ws1.Range("C2").Font.Color = 13391121
$ws1.Range("C2").Font.Name = "Arial"
$ws1.Range("C2").Font.Bold = $false

# D2/E2 are dropped entirely.
$ws1.Range("D2:E2").Hyperlinks.Delete()
$ws1.Range("D2:E2").ClearContents()

# ---------------------------------------------------------------------
# Sheet1 - new data row 3 (mirrors row 2, with a plain-default name cell)
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "sachin"
$ws1.Range("B3").Value = "MyTest"
$ws1.Range("C3").Value = "https://github.com/shivamgupta2607/MyTest"

$ws1.Range("B3").Font.Name = "Arial"

$ws1.Hyperlinks.Add($ws1.Range("C3"), "https://github.com/shivamgupta2607/MyTest") | Out-Null
$ws1.Range("C3").Font.Underline = $true
$ws1.Range("C3").Font.Color = 13391121
$ws1.Range("C3").Font.Name = "Arial"
$ws1.Range("C3").Font.Bold = $false

# ---------------------------------------------------------------------
# Sheet1 - column width (column C gets wider)
# ---------------------------------------------------------------------
$ws1.Columns("C").ColumnWidth = 53.17

# ---------------------------------------------------------------------
# Add Sheet2 after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Test Case Repo"
$ws2.Range("B1").Value = "Test Case Repo URL"
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("B1").Font.Color = 0
$ws2.Range("B1").Font.Name = $null

$ws2.Range("C1:Z1").Font.Bold = $true
$ws2.Range("C1:Z1").Font.Name = "Arial"

$ws2.Range("A2").Value = "interviews"
$ws2.Range("A2").Font.Name = "Arial"

$ws2.Range("B2").Value = "git@repo2.deskera.com:infinity-stones/interviews.git"
$ws2.Range("B2").Font.Color = 0
$ws2.Range("B2").Font.Name = "Arial"

$ws2.Columns("A").ColumnWidth = 14.46
$ws2.Columns("B").ColumnWidth = 44.31
